# "You Are What YOu Eat Final Uploaded"
# Applies the post_info.xlsx edits: updated image filename, trimmed description,
# a real date value for the submission date, tab-ratio / column / row cosmetics,
# and a re-anchored/resized title picture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell content changes
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "./images/tdos_youarewhatyoueat_dairycowsnew.jpeg"
$ws.Range("D11").Value = "We know that the planet is slowly dying, but what you might not know is how much your weekly cheeseburger is contributing to not only greenhouse gas emissions but also water and land depletion."

# D21 used to hold the literal text "2019/11/06"; it becomes a real date value
# formatted as YYYY-MM-DD (serial 43777 == 2019-11-08).
$ws.Range("D21").NumberFormat = "YYYY\-MM\-DD"
$ws.Range("D21").Value = [datetime]"2019-11-08"

# ---------------------------------------------------------------------------
# 2. Column widths (character units = raw xlsx width - 5/6)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.020917678812417
$ws.Columns.Item(2).ColumnWidth = 42.30836707152496
$ws.Columns.Item(3).ColumnWidth = 8.020917678812417
$ws.Columns.Item(4).ColumnWidth = 43.30836707152496
$ws.Columns.Item(5).ColumnWidth = 8.020917678812417
$ws.Columns.Item(6).ColumnWidth = 8.020917678812417
$ws.Columns.Item(7).ColumnWidth = 8.020917678812417
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 13.599865047233466

# ---------------------------------------------------------------------------
# 3. Row heights: several auto-height rows settle back to the 15pt default
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(20).AutoFit()

# ---------------------------------------------------------------------------
# 4. Re-anchor / resize the title picture (in points; 1pt = 12700 EMU)
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 10.998425196850393
$shp.Top = 2.3244094488188978
$shp.Width = 39.68503937007874
$shp.Height = 44.95748031496063

# ---------------------------------------------------------------------------
# 5. View state: selection moves to D10, top-left scrolls back to A1
# ---------------------------------------------------------------------------
$ws.Range("D10").Select()

# ---------------------------------------------------------------------------
# 6. Book-level cosmetic: widen the sheet-tab area (tabRatio 403 -> 600)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.TabRatio = 0.6
